# Update the "想去人数" (want-to-go count) figures for a few events that
# appear on both the "展览" (Exhibition) sheet and the "全部类型" (All types)
# sheet, reflecting newly scraped counts.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1): rows 2, 4, 5 in column F ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 662
$wsExhibit.Range("F4").Value = 1531
$wsExhibit.Range("F5").Value = 708

# --- Sheet "全部类型" (sheet4): rows 2, 4, 6 in column F ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 662
$wsAll.Range("F4").Value = 1531
$wsAll.Range("F6").Value = 708
